# "Generate Report for handback"
#
# This script fills in the handback (target-language received) columns for
# the two localized-file rows on the "zh-cn" and "de-de" sheets, marks the
# status of every tracked file as handed-back, and records the handback
# datetime that was previously a placeholder ("0001-01-01 00:00:00").

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$fileHyperlinkBase = "https://github.com/OpenLocalizationTest/oltest/blob/9c1fb7a2dde0da9b1b51da6b36114bc6d5e5f2e7/e2e/"
$configHyperlink   = "https://github.com/OpenLocalizationTest/oltest/blob/9c1fb7a2dde0da9b1b51da6b36114bc6d5e5f2e7/.localization-config"

$file1 = "5649526d-2658-4a71-ac2a-0d48f2676abd.md"
$file2 = "57291292-9cef-4036-b040-086c0f490e78.md"

$zhXlf1 = "5649526d-2658-4a71-ac2a-0d48f2676abd.08601003cdec70cdc2c6a19e5d7e2a7d5facace1.zh-cn.xlf"
$zhXlf2 = "57291292-9cef-4036-b040-086c0f490e78.b834094a6955464382a6bd05d60654acb0cc6202.zh-cn.xlf"
$zhXlf1Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e3a74691603551c88fcdd08905536df18f4e108/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $zhXlf1
$zhXlf2Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e3a74691603551c88fcdd08905536df18f4e108/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $zhXlf2

$deXlf1 = "5649526d-2658-4a71-ac2a-0d48f2676abd.08601003cdec70cdc2c6a19e5d7e2a7d5facace1.de-de.xlf"
$deXlf2 = "57291292-9cef-4036-b040-086c0f490e78.b834094a6955464382a6bd05d60654acb0cc6202.de-de.xlf"
$deXlf1Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/621b9e3ef433aa44eed6a594a862f3bedfcd6777/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $deXlf1
$deXlf2Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/621b9e3ef433aa44eed6a594a862f3bedfcd6777/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $deXlf2

$zhHandback = "2016-02-15 04:25:43"
$deHandback = "2016-02-15 04:26:08"

# blue underlined "hyperlink" look used throughout the workbook
$hyperlinkColor = 15570276   # 0x6495ED (cornflower blue) as BGR long

function Set-HandbackLink {
    param($ws, $cellRef, $target, $display)
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, [Type]::Missing, [Type]::Missing, $display) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet: status text changes for both rows
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

Set-HandbackLink $zh "E2" ($fileHyperlinkBase + $file1) $file1
Set-HandbackLink $zh "F2" $zhXlf1Target $zhXlf1
$zh.Range("G2").Value = $zhHandback

Set-HandbackLink $zh "E3" ($fileHyperlinkBase + $file2) $file2
Set-HandbackLink $zh "F3" $zhXlf2Target $zhXlf2
$zh.Range("G3").Value = $zhHandback

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

Set-HandbackLink $de "E2" ($fileHyperlinkBase + $file1) $file1
Set-HandbackLink $de "F2" $deXlf1Target $deXlf1
$de.Range("G2").Value = $deHandback

Set-HandbackLink $de "E3" ($fileHyperlinkBase + $file2) $file2
Set-HandbackLink $de "F3" $deXlf2Target $deXlf2
$de.Range("G3").Value = $deHandback

Write-Host "Handback report generated."
